# Daily attendance processing - 2025-12-25 09:35:10
# Swap the order of the two comma-separated "Recorded By" values in column G
# for every row where the cell starts with "dnasr281@gmail.com, ".
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#      "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$needle = "dnasr281@gmail.com, "

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith($needle)) {
        $rest = $val.Substring($needle.Length)
        $newVal = $rest + ", dnasr281@gmail.com"
        $cell.Value2 = $newVal
    }
}
